$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.408.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.567.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.02%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.08%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.78"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.69%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.98"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.27%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.246"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0590"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.79%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0891"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.34%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.792.49"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.98%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.563.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.36%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.65%  "

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.418.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.67"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.49%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.14"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.62%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "226.26"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.76%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0689"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.27%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.88"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.93%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.08"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.95%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.07"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.93"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.66%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.66%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.20%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.69%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.06"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.391.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.26%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.01%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0165"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.36%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.528"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.15%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.90"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.12%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.783"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.44%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.975"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.02%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.46"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.704.01"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.82%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.96%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.64%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0517"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.26%  "
